$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "blue" column (G) with confirmation values for each review row.
$ws.Range("G1").Value = "blue"
$ws.Range("G2").Value = "confirm"
$ws.Range("G3").Value = "confirm"
$ws.Range("G4").Value = "no"
$ws.Range("G5").Value = "confirm"
$ws.Range("G6").Value = "confirm"
$ws.Range("G7").Value = "confirm"
$ws.Range("G8").Value = "confirm"
$ws.Range("G9").Value = "no"
$ws.Range("G10").Value = "yes"

# Match the updated view/selection state from the saved workbook.
$ws.Range("G11").Select()
